$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 891.61536
$ws.Range("I11").Value = 891.61536
$ws.Range("K11").Value = 891.61536
$ws.Range("M11").Value = -751.61536
$ws.Range("H43").Value = 2703.96
$ws.Range("J43").Value = 6214.2856
$ws.Range("L43").Value = 6214.2856
$ws.Range("N43").Value = -6352.2856
$ws.Range("H132").Value = 1886.5128
$ws.Range("I132").Value = 1605.2941
$ws.Range("K132").Value = 4815.8823
$ws.Range("M132").Value = -2285.8823
$ws.Range("H135").Value = 25079.6
$ws.Range("I135").Value = 5000
$ws.Range("J135").Value = 30099.5
$ws.Range("K135").Value = 45000
$ws.Range("L135").Value = 270895.5
$ws.Range("M135").Value = -42465
$ws.Range("N135").Value = -275965.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18527818
$ws.Range("I32").Value = 18527818
$ws.Range("K32").Value = 18527818
$ws.Range("M32").Value = -18527531
$ws.Range("H45").Value = 1618.3636
$ws.Range("I45").Value = 873.0909
$ws.Range("J45").Value = 2363.6365
$ws.Range("K45").Value = 873.0909
$ws.Range("L45").Value = 2363.6365
$ws.Range("M45").Value = -496.0909
$ws.Range("N45").Value = -3117.6365
$ws.Range("H74").Value = 8070926
$ws.Range("I74").Value = 13158729
$ws.Range("J74").Value = 15237.083
$ws.Range("K74").Value = 13158729
$ws.Range("L74").Value = 15237.083
$ws.Range("M74").Value = -13157855
$ws.Range("N74").Value = -16985.083
$ws.Range("H77").Value = 8070926
$ws.Range("I77").Value = 13158729
$ws.Range("J77").Value = 15237.083
$ws.Range("K77").Value = 65793645
$ws.Range("L77").Value = 76185.41500000001
$ws.Range("M77").Value = -65789277
$ws.Range("N77").Value = -84921.41500000001
$ws.Range("H97").Value = 1129
$ws.Range("I97").Value = 1177.375
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1177.375
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -681.375
$ws.Range("N97").Value = -1992
$ws.Range("H122").Value = 1937.9445
$ws.Range("I122").Value = 1592.2667
$ws.Range("K122").Value = 4776.800099999999
$ws.Range("M122").Value = -2326.800099999999
$ws.Range("H125").Value = 67131
$ws.Range("J125").Value = 67131
$ws.Range("L125").Value = 67131
$ws.Range("N125").Value = -76971
$ws.Range("H132").Value = 7259.409
$ws.Range("I132").Value = 3963.0667
$ws.Range("K132").Value = 11889.2001
$ws.Range("M132").Value = -9359.2001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 328.16666
$ws.Range("I22").Value = 328.16666
$ws.Range("K22").Value = 328.16666
$ws.Range("M22").Value = -155.16666
$ws.Range("H105").Value = 1500.8462
$ws.Range("J105").Value = 1321.4
$ws.Range("L105").Value = 1321.4
$ws.Range("N105").Value = -4815.4
$ws.Range("H134").Value = 33831.355
$ws.Range("I134").Value = 1546.5358
$ws.Range("J134").Value = 335156.34
$ws.Range("K134").Value = 4639.607400000001
$ws.Range("L134").Value = 1005469.02
$ws.Range("M134").Value = -2104.607400000001
$ws.Range("N134").Value = -1010539.02

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6842.8
$ws.Range("I7").Value = 466.66666
$ws.Range("K7").Value = 466.66666
$ws.Range("M7").Value = -353.66666
$ws.Range("H58").Value = 2405.0908
$ws.Range("I58").Value = 2120.625
$ws.Range("J58").Value = 3163.6667
$ws.Range("K58").Value = 2120.625
$ws.Range("L58").Value = 3163.6667
$ws.Range("M58").Value = -1917.625
$ws.Range("N58").Value = -3569.6667
$ws.Range("H107").Value = 1329.421
$ws.Range("I107").Value = 890.5833
$ws.Range("K107").Value = 890.5833
$ws.Range("M107").Value = 1029.4167
$ws.Range("H136").Value = 2405.0908
$ws.Range("I136").Value = 2120.625
$ws.Range("J136").Value = 3163.6667
$ws.Range("K136").Value = 6361.875
$ws.Range("L136").Value = 9491.000100000001
$ws.Range("M136").Value = -3811.875
$ws.Range("N136").Value = -14591.0001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 77.166664
$ws.Range("I26").Value = 51.57143
$ws.Range("J26").Value = 113
$ws.Range("K26").Value = 154.71429
$ws.Range("L26").Value = 339
$ws.Range("M26").Value = 133.28571
$ws.Range("N26").Value = -915
$ws.Range("H56").Value = 10900.454
$ws.Range("I56").Value = 10900.454
$ws.Range("K56").Value = 10900.454
$ws.Range("M56").Value = -10370.454
$ws.Range("H109").Value = 599.17645
$ws.Range("I109").Value = 352.4
$ws.Range("J109").Value = 2450
$ws.Range("K109").Value = 1057.2
$ws.Range("L109").Value = 7350
$ws.Range("M109").Value = -17.19999999999982
$ws.Range("N109").Value = -9430

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 260.53845
$ws.Range("I2").Value = 19.428572
$ws.Range("K2").Value = 19.428572
$ws.Range("M2").Value = 93.571428
$ws.Range("H11").Value = 10296660
$ws.Range("I11").Value = 7127029.5
$ws.Range("J11").Value = 17032126
$ws.Range("K11").Value = 7127029.5
$ws.Range("L11").Value = 17032126
$ws.Range("M11").Value = -7126890.5
$ws.Range("N11").Value = -17032404
$ws.Range("H111").Value = 49832
$ws.Range("J111").Value = 49832
$ws.Range("L111").Value = 49832
$ws.Range("N111").Value = -55966
$ws.Range("H132").Value = 33334928
$ws.Range("I132").Value = 33334928
$ws.Range("K132").Value = 100004784
$ws.Range("M132").Value = -100002254

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1891.5625
$ws.Range("I16").Value = 1451.6154
$ws.Range("K16").Value = 1451.6154
$ws.Range("M16").Value = -1281.6154
$ws.Range("H23").Value = 13670.667
$ws.Range("I23").Value = 13670.667
$ws.Range("K23").Value = 13670.667
$ws.Range("M23").Value = -13440.667
$ws.Range("H46").Value = 3449.9285
$ws.Range("I46").Value = 2806.75
$ws.Range("J46").Value = 4307.5
$ws.Range("K46").Value = 2806.75
$ws.Range("L46").Value = 4307.5
$ws.Range("M46").Value = -2618.75
$ws.Range("N46").Value = -4683.5
$ws.Range("H55").Value = 111111550
$ws.Range("I55").Value = 125000450
$ws.Range("K55").Value = 125000450
$ws.Range("M55").Value = -125000277
$ws.Range("H61").Value = 700
$ws.Range("I61").Value = 700
$ws.Range("K61").Value = 700
$ws.Range("M61").Value = -498
$ws.Range("H113").Value = 700
$ws.Range("I113").Value = 700
$ws.Range("K113").Value = 700
$ws.Range("M113").Value = 1470
$ws.Range("H127").Value = 125306.664
$ws.Range("J127").Value = 125306.664
$ws.Range("L127").Value = 125306.664
$ws.Range("N127").Value = -135226.664

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1083.8
$ws.Range("I23").Value = 140
$ws.Range("J23").Value = 2499.5
$ws.Range("K23").Value = 140
$ws.Range("L23").Value = 2499.5
$ws.Range("M23").Value = 89
$ws.Range("N23").Value = -2957.5
$ws.Range("H92").Value = 63030
$ws.Range("J92").Value = 63030
$ws.Range("L92").Value = 63030
$ws.Range("N92").Value = -68022
$ws.Range("H94").Value = 68325
$ws.Range("J94").Value = 68325
$ws.Range("L94").Value = 68325
$ws.Range("N94").Value = -70127
